$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (nombre_aides, montant_total)  -- stored as text, same as original cells
$updates = @{
    31  = @("441", "1300749.11")
    33  = @("819", "5354648.45")
    35  = @("548", "2974242.32")
    38  = @("591", "1547372.32")
    39  = @("293", "1590068.04")
    40  = @("275", "928520.72")
    50  = @("996", "6372799.81")
    58  = @("6950", "35597073.29")
    60  = @("6812", "29488918.34")
    62  = @("137", "692707.46")
    100 = @("1376", "3483754.28")
    103 = @("1563", "7975645.55")
    105 = @("1509", "7035240.06")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $vals[0]
    $cellC.Style = "Normal"

    $cellD = $ws.Cells.Item($row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $vals[1]
    $cellD.Style = "Normal"
}
